$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "27.696.58"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.633.34"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'212.06"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D8").Value = "'23.20"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "'0.0862"
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("D12").Value = "1.865.45"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "1.638.68"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "'65.18"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "27.667.83"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'229.51"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "'7.56"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +4.67%  "
$ws.Range("D23").Value = "'4.34"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "'149.03"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").Value = "'15.57"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "'3.28"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").Value = "1.471.81"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("D37").Value = "'0.937"
$ws.Range("E37").Value = "  +3.41%  "
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").Value = "'0.876"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "'67.78"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("E45").Value = "  -4.38%  "
$ws.Range("D46").Value = "1.774.77"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").Value = "'87.64"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("D50").Value = "'0.0993"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -1.46%  "
